# Apply the table-style change described by the commit:
#   slide 5's table switches from table style
#   {99846C70-C94D-4242-AAD0-E9CE491DF5EB} ("Table_0", the custom pink style
#   declared in ppt/tableStyles.xml) to the built-in table style
#   {539C43E2-4BF7-4072-9FDB-7A40AFBCB7DF}.
#
# This mirrors a user selecting a different style from the Table Styles
# gallery (Table Tools > Design) for the single table on slide 5.

$p = $ppt.ActivePresentation

# The table lives on slide 5 ("B1- TYPES OF FINANCIAL DOCUMENTS"), as the
# second shape (a graphicFrame hosting the a:tbl).
$slide = $p.Slides.Item(5)
$tableShape = $slide.Shapes.Item(2)

if ($tableShape.HasTable) {
    $table = $tableShape.Table
    $table.ApplyStyle("{539C43E2-4BF7-4072-9FDB-7A40AFBCB7DF}")
}
